$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New incident/ticket rows to append after the existing data (rows 182-186).
# Columns: A=Fecha, B=Hora, C=WC47 NACP, D=WC48 P5F, E=WC49 P5H, F=WV50 FILTER,
#          G=SPL, H=Hora de Reparacion, I=Tiempo de Reparacion
# Row 185 has no repair time / duration recorded (columns H and I stay empty).
$newRows = @(
    @{ Row = 182; A = "2024-05-22"; B = "12:02:22"; C = "Fallo en elevador_2";  D = "-"; E = "-"; F = "-"; G = "-"; H = "12:02:23"; I = "0:00:01" },
    @{ Row = 183; A = "2024-05-22"; B = "12:07:33"; C = "Fallo en paletizador"; D = "-"; E = "-"; F = "-"; G = "-"; H = "12:07:40"; I = "0:00:07" },
    @{ Row = 184; A = "2024-05-22"; B = "12:07:37"; C = "No atornilla clips";   D = "-"; E = "-"; F = "-"; G = "-"; H = "12:07:41"; I = "0:00:04" },
    @{ Row = 185; A = "2024-05-22"; B = "12:07:59"; C = "Fallo en elevador";    D = "-"; E = "-"; F = "-"; G = "-" },
    @{ Row = 186; A = "2024-05-22"; B = "12:08:28"; C = "Fallo en paletizador"; D = "-"; E = "-"; F = "-"; G = "-"; H = "12:08:31"; I = "0:00:03" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A holds a date-looking string ("2024-05-22"). Left as a normal
    # value it would be auto-parsed into a numeric date serial, but the
    # source data stores it as plain text, so force text interpretation by
    # switching the cell to a text format before assigning the value, then
    # switch the cell's style back to Normal so the cell keeps the text
    # value without leaving a visible custom number format behind.
    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $r.A
    $aCell.Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G

    if ($r.ContainsKey("H")) {
        $ws.Cells.Item($rowNum, 8).Value = $r.H
    }
    if ($r.ContainsKey("I")) {
        $ws.Cells.Item($rowNum, 9).Value = $r.I
    }
}
